$wb = $excel.ActiveWorkbook

# Thing
$wb.Worksheets.Item("Thing").Range("A1").Value = "identifier"

# Intangible
$wb.Worksheets.Item("Intangible").Range("A1").Value = "identifier"

# CreativeWork
$wb.Worksheets.Item("CreativeWork").Range("A1").Value = "identifier"

# DefinedTerm
$wb.Worksheets.Item("DefinedTerm").Range("C1").Value = "identifier"

# DefinedTermSet
$wb.Worksheets.Item("DefinedTermSet").Range("B1").Value = "identifier"
